# Apply the commit:
#  1. "bernoulli" -> sheet "Coefficients": row 4 (College Graduates) switched from the
#     repeated placeholder formulas to real bootstrap/bernoulli coefficient values,
#     highlighted in yellow.
#  2. "take out the variance of labor income and unemployment income" -> sheet
#     "Variance": split the single "Panel A - Labor Income Only" table into two
#     side-by-side panels - "Labor Income Only" (existing columns) and a new
#     "Labor Income Plus Unemployment Income" panel (columns E:G), with the
#     differing/added statistics highlighted in yellow.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Coefficients"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Coefficients")

# Row 4 ("College Graduates") gets new literal values (no longer formulas),
# highlighted with a yellow fill.
$ws1.Range("B4").Value = -51819.03
$ws1.Range("C4").Value = 3553.459
$ws1.Range("D4").Value = 0.9337594
$ws1.Range("E4").Value = -0.5499252
$ws1.Range("B4:E4").Interior.Color = 65535

# ---------------------------------------------------------------------------
# Sheet 2: "Variance"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Variance")

# Mirror columns B:D (widths + number formats) into new columns E:G so the
# second panel looks just like the first.
$ws2.Columns("E").ColumnWidth = $ws2.Columns("B").ColumnWidth
$ws2.Columns("F").ColumnWidth = $ws2.Columns("C").ColumnWidth
$ws2.Columns("G").ColumnWidth = $ws2.Columns("D").ColumnWidth

# Relabel the first panel title/subtitle.
$ws2.Range("A2").Value = "Labor Income Only"

# Second panel title (row 1, merged E1:G1) + subtitle (row 2, merged E2:G2),
# copying the formatting of the first panel's title/subtitle cells.
$ws2.Range("A1:D1").Copy() | Out-Null
$ws2.Range("E1:G1").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2:D2").Copy() | Out-Null
$ws2.Range("E2:G2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws2.Range("E1:G1").Merge() | Out-Null
$ws2.Range("E2:G2").Merge() | Out-Null

$ws2.Range("E1").Value = " Table 2 - Variance Decomposition"
$ws2.Range("E2").Value = "Labor Income Plus Unemployment Income"

# Column headers for the second panel (row 3), copying the first panel's header format.
$ws2.Range("B3:D3").Copy() | Out-Null
$ws2.Range("E3:G3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws2.Range("E3").Value = "No High School"
$ws2.Range("F3").Value = "High School Graduates"
$ws2.Range("G3").Value = "College Graduates"

# Data rows 4-7, copying the first panel's number format into the second panel.
$ws2.Range("B4:D4").Copy() | Out-Null
$ws2.Range("E4:G4").PasteSpecial(-4122) | Out-Null
$ws2.Range("B5:D5").Copy() | Out-Null
$ws2.Range("E5:G5").PasteSpecial(-4122) | Out-Null
$ws2.Range("B6:D6").Copy() | Out-Null
$ws2.Range("E6:G6").PasteSpecial(-4122) | Out-Null
$ws2.Range("B7:D7").Copy() | Out-Null
$ws2.Range("E7:G7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Panel 1 (labor income only) updated values.
$ws2.Range("B4").Value = 0.15716233645501712
$ws2.Range("C4").Value = 0.13711309200802088
$ws2.Range("D4").Value = 0.14387494569938158

$ws2.Range("B6").Value = 0.3391902121229326
$ws2.Range("C6").Value = 0.23377339455121918
$ws2.Range("D6").Value = 0.22956480566497992

# Panel 2 (labor income + unemployment income) values.
$ws2.Range("E4").Value = 0.16911534525287764
$ws2.Range("F4").Value = 0.1374772708486752
$ws2.Range("G4").Value = 0.14525839046333949

$ws2.Range("E5").Value = -0.0029
$ws2.Range("F5").Value = -0.0017
$ws2.Range("G5").Value = -0.0017

$ws2.Range("E6").Value = 0.38529209698616967
$ws2.Range("F6").Value = 0.28178005607210743
$ws2.Range("G6").Value = 0.26258332011001767

$ws2.Range("E7").Value = -0.0182
$ws2.Range("F7").Value = -0.0182
$ws2.Range("G7").Value = -0.0107

# Highlight the cells that now differ between the two panels / are new.
$ws2.Range("D4").Interior.Color = 65535
$ws2.Range("G4").Interior.Color = 65535
$ws2.Range("D6").Interior.Color = 65535
$ws2.Range("G6").Interior.Color = 65535
$ws2.Range("E7:G7").Borders.Item(9).LineStyle = 1

# ---------------------------------------------------------------------------
# View state: zoom + selection on both sheets, "Coefficients" left as the
# selected/active tab.
# ---------------------------------------------------------------------------
$ws2.Select() | Out-Null
$excel.ActiveWindow.Zoom = 150
$ws2.Range("C9").Select() | Out-Null

$ws1.Select() | Out-Null
$excel.ActiveWindow.Zoom = 187
$ws1.Range("D12").Select() | Out-Null
